# slight edit to data request
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "X" marks in column F for rows 3-6 (Data Set 1 section)
$ws.Range("F3").Value = "X"
$ws.Range("F4").Value = "X"
$ws.Range("F5").Value = "X"
$ws.Range("F6").Value = "X"

# Add "X" mark in column B for row 17 (Academic data section)
$ws.Range("B17").Value = "X"

# Update the selected cell to C11
$ws.Range("C11").Select()
